# "Caso de busca de produto negativo - Sucesso"
# Update the test-user row (usuario/email) used by the negative product
# search test case, and leave the selection where the tester left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# usuario / email columns (A2:B2) for the negative-search test case.
$ws.Range("A2").Value = "Braulioteste4"
$ws.Range("B2").Value = "teste4@teste.com.br"

# Make sure we're on the right sheet, then move the selection to D7
# (matches the saved cursor position in the workbook).
$ws.Activate()
$ws.Range("D7").Select()
